$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.560068368911743
$ws.Range("B1").Value = 2.414801120758057
$ws.Range("C1").Value = 1.773718118667603
$ws.Range("D1").Value = 1.617916584014893
$ws.Range("E1").Value = 1.437620162963867
